# derivative/manifest_expected.xlsx - test_annotate_plot
# Reorders the "annotated plot" file-group rows so that the order becomes
# proximal, transverse, distal (was distal, proximal, transverse), and
# refreshes the column layout / selection to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8-9: now "proximal" (was "distal") ---
$ws.Range("A8").Value = "stim_proximal-colon_manometry.csv"
$ws.Range("H8").Value = "stim_proximal-colon_manometry.jpg"
$ws.Range("A9").Value = "stim_proximal-colon_manometry.jpg"
$ws.Range("I9").Value = "stim_proximal-colon_manometry.csv"

# --- Row 10-11: now "transverse" (was "proximal") ---
$ws.Range("A10").Value = "stim_transverse-colon_manometry.csv"
$ws.Range("H10").Value = "stim_transverse-colon_manometry.jpg"
$ws.Range("A11").Value = "stim_transverse-colon_manometry.jpg"
$ws.Range("I11").Value = "stim_transverse-colon_manometry.csv"

# --- Row 12-13: now "distal" (was "transverse") ---
$ws.Range("A12").Value = "stim_distal-colon_manometry.csv"
$ws.Range("H12").Value = "stim_distal-colon_manometry.jpg"
$ws.Range("A13").Value = "stim_distal-colon_manometry.jpg"
$ws.Range("I13").Value = "stim_distal-colon_manometry.csv"

# --- Column layout: best-fit widths, with the two long JSON/description
#     columns (C, E) hidden, matching the saved workbook view ---
$ws.Columns.Item(1).ColumnWidth = 42.6640625
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 134.5
$ws.Columns.Item(4).ColumnWidth = 7.5
$ws.Columns.Item(5).ColumnWidth = 162.6640625
$ws.Columns.Item(6).ColumnWidth = 26.1640625
$ws.Columns.Item(7).ColumnWidth = 48.33203125
$ws.Columns.Item(8).ColumnWidth = 30.83203125
$ws.Columns.Item(9).ColumnWidth = 30.83203125

$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(5).Hidden = $true

# --- Selection moved to A31 (reflects the last-saved cursor position) ---
[void]$ws.Range("A31").Select()
